$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet1: ip_address_list
# ======================================================================
$ws1 = $wb.Worksheets.Item(1)

# Drop the "440_Austin" row (old row 1) entirely; everything below shifts up.
$ws1.Rows.Item(1).Delete()

# After that shift, old rows 7 & 8 (515_ZF Stara Boleslav, duplicate 503_Witte)
# are now rows 6 & 7 - drop both.
$ws1.Rows.Item(6).Delete()
$ws1.Rows.Item(6).Delete()

# Old row 4 (514_Teleflex) is now row 3 - drop it too (it lives on only in
# ip_adress_fav_list / disk_list).
$ws1.Rows.Item(3).Delete()

# Remaining rows are now: 497_Edcha, 503_Witte, 518_Valeo, 527_Teijin.
# Insert a fresh row for "518_Valeo II" right after 518_Valeo (row 3).
$ws1.Rows.Item(4).Insert()
$ws1.Range("A4").Value = "518_Valeo II"
$ws1.Range("B4").Value = "192.168.1.243"
$ws1.Range("C4").Value = "255.255.255.0"
$ws1.Range("E4").Value = $true

# 527_Teijin (now row 5) gets an extra "OP:" line in its notes, and its
# favourite flag becomes a plain 0 instead of boolean FALSE.
$ws1.Range("D5").Value = "XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28.117"
$ws1.Range("E5").Value = 0

# Append the brand-new "511_Teleflex" row at the bottom.
$ws1.Range("A6").Value = "511_Teleflex"
$ws1.Range("B6").Value = "192.168.1.242"
$ws1.Range("C6").Value = "255.255.255.0"
$ws1.Range("D6").Value = "Teleflex "
$ws1.Range("E6").Value = 0

# ======================================================================
# Sheet2: ip_adress_fav_list
# ======================================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1:E3").ClearContents()
$ws2.Range("A1").Value = "518_Valeo II"
$ws2.Range("B1").Value = "192.168.1.243"
$ws2.Range("C1").Value = "255.255.255.0"
$ws2.Range("E1").Value = 1

# ======================================================================
# Sheet3: disk_list
# ======================================================================
$ws3 = $wb.Worksheets.Item(3)

# Row 1 (514_Teleflex) is untouched. Drop "Domaci Nas" and "518_Valeo II"
# (old rows 3 & 4); this leaves 514_Teleflex, 515_ZF, 518_Valeo, 474_B Austin.
$ws3.Rows.Item(3).Delete()
$ws3.Rows.Item(3).Delete()

# Move 515_ZF (now row 2) down below 474_B Austin (row 4) -> becomes row 4.
$ws3.Range("A2:F2").Cut($ws3.Range("A5:F5"))
$ws3.Rows.Item(2).Delete()

# ======================================================================
# Sheet4: Settings
# ======================================================================
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B6").Value = 0
$ws4.Range("B8").Value = 1
$ws4.Range("A10").Value = "statusy odvolby dotazov�n� p�i maz�n�"
$ws4.Range("B10").Value = 110

# Settings becomes the active tab/sheet (matches activeTab=3 in the saved file).
$ws4.Activate()

# ======================================================================
# New hidden sheet: projects_bin2 (recycle bin for deleted projects)
# ======================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$binSheet = $wb.Worksheets.Add($null, $lastSheet)
$binSheet.Name = "projects_bin2"
$binSheet.Visible = $false
